$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "HQ box" note from B4 (Burning Hands) ...
$ws.Range("B4").ClearContents()

# ... and move it down to B14 (Sleep)
$ws.Range("B14").Value = "HQ box"

# Add a new location note for C6 (Chill Touch row)
$ws.Range("C6").Value = "Skyflow -> 4F -> Main Hall"

# Add a new location note for B7 (Chromatic Orb row)
$ws.Range("B7").Value = "Human Resource basement"

# Update the selected cell to reflect where work left off
$ws.Range("F8").Select()
